$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1: copy formatting (style) from the existing H1 header
# cell, then set the header text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-71: each entry is (row, I value, J value)
$ijData = @(
    @(2,8,8),
    @(3,7,7),
    @(4,8,9),
    @(5,6,7),
    @(6,7,7),
    @(7,9,9),
    @(8,7,7),
    @(9,7,7),
    @(10,7,7),
    @(11,7,7),
    @(12,6,7),
    @(13,9,9),
    @(14,7,7),
    @(15,7,7),
    @(16,6,7),
    @(17,9,9),
    @(18,7,7),
    @(19,8,8),
    @(20,5,5),
    @(21,6,7),
    @(22,5,6),
    @(23,7,7),
    @(24,4,5),
    @(25,5,5),
    @(26,7,7),
    @(27,8,8),
    @(28,6,7),
    @(29,9,9),
    @(30,7,8),
    @(31,7,7),
    @(32,8,8),
    @(33,6,6),
    @(34,5,6),
    @(35,7,8),
    @(36,8,8),
    @(37,7,8),
    @(38,7,7),
    @(39,7,7),
    @(40,6,7),
    @(41,8,8),
    @(42,6,8),
    @(43,7,8),
    @(44,6,6),
    @(45,7,9),
    @(46,7,8),
    @(47,8,8),
    @(48,7,7),
    @(49,7,7),
    @(50,8,8),
    @(51,6,7),
    @(52,8,8),
    @(53,7,7),
    @(54,7,7),
    @(55,8,9),
    @(56,8,8),
    @(57,8,8),
    @(58,7,8),
    @(59,6,7),
    @(60,6,6),
    @(61,9,9),
    @(62,8,9),
    @(63,9,9),
    @(64,9,9),
    @(65,9,9),
    @(66,9,9),
    @(67,7,7),
    @(68,7,7),
    @(69,5,5),
    @(70,6,6),
    @(71,6,6)
)

foreach ($entry in $ijData) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
